# "Steal or Deal" section - update the round thresholds / gold rewards.
$d = $word.ActiveDocument

# Max gold the player can end with: 80 -> 100
$d.Content.Find.Execute("In this game you can end with 80 gold.", $true, $false, $false, $false, $false, $true, 1, $false, "In this game you can end with 100 gold.", 2)

# Starting gold line.
$d.Content.Find.Execute("You start with 5 gold", $true, $false, $false, $false, $false, $true, 1, $false, "Starts at 0 gold", 2)

# Round thresholds, shifted down a bracket and re-priced.
$d.Content.Find.Execute("Between 10-20 = 10 gold", $true, $false, $false, $false, $false, $true, 1, $false, "Between 1-9 = 5 gold", 2)
$d.Content.Find.Execute("Between 20-30 = 15 gold", $true, $false, $false, $false, $false, $true, 1, $false, "Between 10-19 = 10 gold", 2)
$d.Content.Find.Execute("Between 30-40 = 20 gold", $true, $false, $false, $false, $false, $true, 1, $false, "Between 20-29 = 15 gold", 2)
$d.Content.Find.Execute("Between 40-50 = 25 gold", $true, $false, $false, $false, $false, $true, 1, $false, "Between 30-39 = 20 gold", 2)
$d.Content.Find.Execute("Between 50-60 = 35 gold", $true, $false, $false, $false, $false, $true, 1, $false, "Between 40-49 = 25 gold", 2)
$d.Content.Find.Execute("Between 60-70 = 45 gold", $true, $false, $false, $false, $false, $true, 1, $false, "Between 50-59 = 35 gold", 2)
$d.Content.Find.Execute("Between 70-80 = 55 gold", $true, $false, $false, $false, $false, $true, 1, $false, "Between 60-69 = 45 gold", 2)
$d.Content.Find.Execute("Between 80-90 = 65 gold", $true, $false, $false, $false, $false, $true, 1, $false, "Between 70-79 = 60 gold", 2)
$d.Content.Find.Execute("Between 90-100 = 80 gold", $true, $false, $false, $false, $false, $true, 1, $false, "Between 80-89 = 80 gold", 2)

# Add the new top bracket (90-100) as a new list item right after "Between 80-89 = 80 gold".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $ptext = $para.Range.Text.TrimEnd([char]13)
    if ($ptext -eq "Between 80-89 = 80 gold") {
        $para.Range.InsertParagraphAfter()
        $d.Paragraphs.Item($i + 1).Range.Text = "Between 90-100 = 100 gold"
        break
    }
}
